$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "kepsek"
$ws.Range("C3").Value = "Kepala Sekolah"
